$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 325, shifting the existing rows 325..382 down to 326..383.
$ws.Rows.Item(325).Insert()

# Populate the newly inserted row 325 with the new record's data.
$ws.Cells.Item(325, 1).Value = 7
$ws.Cells.Item(325, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(325, 3).Value = "Ñuble"
$ws.Cells.Item(325, 4).Value = 45209
$ws.Cells.Item(325, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(325, 5).Value = 16
$ws.Cells.Item(325, 6).Value = 100112043
$ws.Cells.Item(325, 7).Value = "Pepino ensalada"
$ws.Cells.Item(325, 8).Value = "Sin especificar"
$ws.Cells.Item(325, 9).Value = "Primera"
$ws.Cells.Item(325, 10).Value = 60
$ws.Cells.Item(325, 11).Value = 14000
$ws.Cells.Item(325, 12).Value = 14000
$ws.Cells.Item(325, 13).Value = 14000
$ws.Cells.Item(325, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(325, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(325, 16).Value = 233
$ws.Cells.Item(325, 17).Value = 60
$ws.Cells.Item(325, 18).Value = "Hortaliza"

$wb.Save()
